# Apply odds/score updates from the 2025-04-02 FlashScore refresh
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("AE2").Value = 21
$ws.Range("AG2").Value = 21
$ws.Range("AH2").Value = 81
$ws.Range("G2").Value = 1.4
$ws.Range("W2").Value = 11
$ws.Range("L3").Value = 1.2
$ws.Range("M3").Value = 4.5
$ws.Range("AD4").Value = 451
$ws.Range("R4").Value = 1.95
$ws.Range("S4").Value = 1.8
$ws.Range("AE5").Value = 19
$ws.Range("AB7").Value = 23
$ws.Range("AC7").Value = 81
$ws.Range("AF7").Value = 41
$ws.Range("G7").Value = 1.48
$ws.Range("I7").Value = 6.5
$ws.Range("N7").Value = 1.95
$ws.Range("O7").Value = 1.95
$ws.Range("P7").Value = 1.4
$ws.Range("Q7").Value = 2.75
$ws.Range("R7").Value = 2.05
$ws.Range("S7").Value = 1.7
$ws.Range("T7").Value = 7
$ws.Range("U7").Value = 7
$ws.Range("Y7").Value = 34
$ws.Range("AE8").Value = 15
$ws.Range("K8").Value = 13
$ws.Range("N8").Value = 1.83
$ws.Range("O8").Value = 2.03
$ws.Range("AA10").Value = 11
$ws.Range("AE10").Value = 23
$ws.Range("AG10").Value = 23
$ws.Range("AH10").Value = 101
$ws.Range("G10").Value = 1.3
$ws.Range("H10").Value = 5.5
$ws.Range("I10").Value = 8.5
$ws.Range("J10").Value = 1.02
$ws.Range("K10").Value = 19
$ws.Range("L10").Value = 1.14
$ws.Range("M10").Value = 5.5
$ws.Range("N10").Value = 1.48
$ws.Range("O10").Value = 2.6
$ws.Range("P10").Value = 1.25
$ws.Range("Q10").Value = 3.75
$ws.Range("R10").Value = 1.8
$ws.Range("S10").Value = 1.95
$ws.Range("T10").Value = 9
$ws.Range("U10").Value = 7.5
$ws.Range("V10").Value = 9
$ws.Range("Z10").Value = 17
$ws.Range("AB11").Value = 17
$ws.Range("AD11").Value = 401
$ws.Range("AE11").Value = 6.5
$ws.Range("AF11").Value = 9.5
$ws.Range("AG11").Value = 9
$ws.Range("AH11").Value = 19
$ws.Range("AI11").Value = 19
$ws.Range("AJ11").Value = 34
$ws.Range("G11").Value = 3.4
$ws.Range("I11").Value = 2.15
$ws.Range("J11").Value = 1.07
$ws.Range("K11").Value = 8.5
$ws.Range("N11").Value = 2.2
$ws.Range("O11").Value = 1.65
$ws.Range("P11").Value = 1.5
$ws.Range("Q11").Value = 2.5
$ws.Range("R11").Value = 1.95
$ws.Range("S11").Value = 1.8
$ws.Range("U11").Value = 17
$ws.Range("V11").Value = 12
$ws.Range("W11").Value = 41
$ws.Range("X11").Value = 29
$ws.Range("Y11").Value = 41
$ws.Range("Z11").Value = 8.5
$ws.Range("N12").Value = 1.67
$ws.Range("O12").Value = 2.15
$ws.Range("AA13").Value = 6.5
$ws.Range("AF13").Value = 12
$ws.Range("AG13").Value = 9
$ws.Range("AH13").Value = 21
$ws.Range("AI13").Value = 17
$ws.Range("AJ13").Value = 23
$ws.Range("G13").Value = 3.2
$ws.Range("H13").Value = 3.3
$ws.Range("I13").Value = 2.25
$ws.Range("J13").Value = 1.04
$ws.Range("K13").Value = 13
$ws.Range("L13").Value = 1.22
$ws.Range("M13").Value = 4
$ws.Range("N13").Value = 1.8
$ws.Range("P13").Value = 1.36
$ws.Range("Q13").Value = 3
$ws.Range("R13").Value = 1.62
$ws.Range("S13").Value = 2.2
$ws.Range("T13").Value = 12
$ws.Range("U13").Value = 17
$ws.Range("V13").Value = 12
$ws.Range("W13").Value = 34
$ws.Range("X13").Value = 23
$ws.Range("Z13").Value = 12
$ws.Range("AB14").Value = 19
$ws.Range("AE14").Value = 51
$ws.Range("AG14").Value = 29
$ws.Range("AH14").Value = 101
$ws.Range("AJ14").Value = 34
$ws.Range("G14").Value = 1.25
$ws.Range("H14").Value = 7.5
$ws.Range("I14").Value = 8
$ws.Range("M14").Value = 17
$ws.Range("V14").Value = 13
$ws.Range("AB16").Value = 19
$ws.Range("AC16").Value = 81
$ws.Range("AE16").Value = 7
$ws.Range("AG16").Value = 12
$ws.Range("AJ16").Value = 41
$ws.Range("H16").Value = 3.1
$ws.Range("I16").Value = 3.1
$ws.Range("J16").Value = 1.11
$ws.Range("K16").Value = 6.5
$ws.Range("L16").Value = 1.53
$ws.Range("M16").Value = 2.38
$ws.Range("N16").Value = 2.7
$ws.Range("O16").Value = 1.44
$ws.Range("P16").Value = 1.62
$ws.Range("Q16").Value = 2.2
$ws.Range("R16").Value = 2.2
$ws.Range("S16").Value = 1.62
$ws.Range("T16").Value = 6
$ws.Range("V16").Value = 10
$ws.Range("X16").Value = 23
$ws.Range("Z16").Value = 6.5
$ws.Range("AH17").Value = 34
$ws.Range("N17").Value = 2
$ws.Range("O17").Value = 1.85
$ws.Range("AE18").Value = 6
$ws.Range("AF18").Value = 10
$ws.Range("AG18").Value = 10
$ws.Range("AH18").Value = 23
$ws.Range("AI18").Value = 23
$ws.Range("G18").Value = 3.2
$ws.Range("I18").Value = 2.4
$ws.Range("J18").Value = 1.1
$ws.Range("K18").Value = 7
$ws.Range("R18").Value = 2.1
$ws.Range("S18").Value = 1.67
$ws.Range("T18").Value = 7.5
$ws.Range("U18").Value = 15
$ws.Range("W18").Value = 34
$ws.Range("AF19").Value = 23
$ws.Range("AG19").Value = 15
$ws.Range("AH19").Value = 51
$ws.Range("G19").Value = 1.75
$ws.Range("H19").Value = 3.75
$ws.Range("I19").Value = 4.5
$ws.Range("N19").Value = 1.95
$ws.Range("O19").Value = 1.9
$ws.Range("T19").Value = 7
$ws.Range("W19").Value = 13
$ws.Range("AG20").Value = 10
$ws.Range("AH20").Value = 26
$ws.Range("G20").Value = 2.75
$ws.Range("I20").Value = 2.7
$ws.Range("J20").Value = 1.08
$ws.Range("K20").Value = 8
$ws.Range("N20").Value = 2.25
$ws.Range("O20").Value = 1.62
$ws.Range("R20").Value = 1.8
$ws.Range("S20").Value = 1.95
$ws.Range("V20").Value = 11
$ws.Range("W20").Value = 29
$ws.Range("AA21").Value = 6
$ws.Range("AD21").Value = 301
$ws.Range("AE21").Value = 8.5
$ws.Range("AF21").Value = 15
$ws.Range("AH21").Value = 34
$ws.Range("AI21").Value = 26
$ws.Range("AJ21").Value = 34
$ws.Range("G21").Value = 2.35
$ws.Range("H21").Value = 3.3
$ws.Range("I21").Value = 3.1
$ws.Range("L21").Value = 1.33
$ws.Range("M21").Value = 3.25
$ws.Range("N21").Value = 2.1
$ws.Range("O21").Value = 1.7
$ws.Range("P21").Value = 1.44
$ws.Range("Q21").Value = 2.63
$ws.Range("R21").Value = 1.91
$ws.Range("S21").Value = 1.91
$ws.Range("T21").Value = 7.5
$ws.Range("U21").Value = 11
$ws.Range("Z21").Value = 9
$ws.Range("AB22").Value = 15
$ws.Range("AE22").Value = 9.5
$ws.Range("AF22").Value = 17
$ws.Range("AG22").Value = 12
$ws.Range("H22").Value = 3.6
$ws.Range("K22").Value = 9.5
$ws.Range("U22").Value = 9.5
$ws.Range("W22").Value = 19
$ws.Range("AD23").Value = 501
$ws.Range("AB25").Value = 19
$ws.Range("AC25").Value = 67
$ws.Range("AD25").Value = 501
$ws.Range("AE25").Value = 9
$ws.Range("G25").Value = 2.05
$ws.Range("H25").Value = 3.2
$ws.Range("I25").Value = 3.9
$ws.Range("J25").Value = 1.08
$ws.Range("K25").Value = 7.5
$ws.Range("L25").Value = 1.44
$ws.Range("M25").Value = 2.63
$ws.Range("N25").Value = 2.35
$ws.Range("O25").Value = 1.57
$ws.Range("P25").Value = 1.53
$ws.Range("Q25").Value = 2.38
$ws.Range("R25").Value = 2.05
$ws.Range("S25").Value = 1.7
$ws.Range("T25").Value = 6
$ws.Range("U25").Value = 8.5
$ws.Range("V25").Value = 9.5
$ws.Range("X25").Value = 19
$ws.Range("Y25").Value = 34
$ws.Range("Z25").Value = 7.5
$ws.Range("AB26").Value = 26
$ws.Range("AC26").Value = 67
$ws.Range("AD26").Value = 401
$ws.Range("AG26").Value = 34
$ws.Range("AH26").Value = 201
$ws.Range("AI26").Value = 81
$ws.Range("AJ26").Value = 67
$ws.Range("G26").Value = 1.18
$ws.Range("H26").Value = 7.5
$ws.Range("I26").Value = 13
$ws.Range("R26").Value = 2
$ws.Range("S26").Value = 1.75
$ws.Range("T26").Value = 9.5
$ws.Range("U26").Value = 7
$ws.Range("V26").Value = 10
$ws.Range("W26").Value = 7
$ws.Range("Y26").Value = 29
$ws.Range("Z26").Value = 21
$ws.Range("AE27").Value = 6.5
$ws.Range("AH27").Value = 15
$ws.Range("AI27").Value = 17
$ws.Range("G27").Value = 3.75
$ws.Range("H27").Value = 3.3
$ws.Range("I27").Value = 1.9
$ws.Range("J27").Value = 1.06
$ws.Range("K27").Value = 8
$ws.Range("L27").Value = 1.36
$ws.Range("M27").Value = 3
$ws.Range("N27").Value = 2.2
$ws.Range("O27").Value = 1.65
$ws.Range("P27").Value = 1.5
$ws.Range("Q27").Value = 2.5
$ws.Range("R27").Value = 2
$ws.Range("S27").Value = 1.73
$ws.Range("T27").Value = 10
$ws.Range("Z27").Value = 8
$ws.Range("AA28").Value = 7
$ws.Range("G28").Value = 2.25
$ws.Range("H28").Value = 3.7
$ws.Range("I28").Value = 2.7
$ws.Range("N28").Value = 1.8
$ws.Range("O28").Value = 2
$ws.Range("X28").Value = 19
$ws.Range("N30").Value = 1.93
$ws.Range("O30").Value = 1.93

# These two odds are no longer available and become blank
$ws.Range("J14").ClearContents()
$ws.Range("K14").ClearContents()
